# Oprava excel tabulek - hodnoty C/G sloupcu pro radky 4-8 se deli 10
# (sjednoceni jednotek) a prida se novy sloupec O s prepoctem pro mBar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Uprava existujicich hodnot (deleni 10) ---
$ws.Range("C4").Value = 2.386517
$ws.Range("G4").Value = 0.38744000000000001

$ws.Range("C5").Value = 2.2905500000000001
$ws.Range("G5").Value = 0.36835166699999999

$ws.Range("C6").Value = 1.83595
$ws.Range("G6").Value = 0.27230833332999999

$ws.Range("C8").Value = 1.6173500000000001
$ws.Range("G8").Value = 0.26400000000000001

# --- Novy sloupec O: prepocet Sklon/Posun mBar pro kazdou hodnotu M ---
$ws.Range("O4").Formula = "=`$G`$4*M4+`$H`$4"
$ws.Range("O5:O33").Formula = "=`$G`$4*M5+`$H`$4"

# --- Vyber aktivni bunky jako v ulozenem sesitu ---
$ws.Range("R20").Select() | Out-Null
